$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.458696246147156
$ws.Range("B1").Value = 3.381775140762329
$ws.Range("C1").Value = 4.236068725585938
$ws.Range("D1").Value = 2.183442831039429
$ws.Range("E1").Value = 0.7384341359138489
